$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mif"
$ws.Range("C2").Value = "Cxcr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.89761366666667
$ws.Range("H2").Value = 38.692841
$ws.Range("I2").Value = 0.1295258291743358
$ws.Range("J2").Value = 0.1295258291743358
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 161.9042856666667
$ws.Range("N2").Value = 485.712857
$ws.Range("O2").Value = 0.3829374249648381
$ws.Range("P2").Value = 0.3829374249648381
$ws.Range("Q2").Value = 2088.178927506304
$ws.Range("R2").Value = 18793.61034755674
$ws.Range("S2").Value = 0.04960028749045565
$ws.Range("T2").Value = 0.04960028749045565

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mif"
$ws.Range("C3").Value = "Cxcr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.89761366666667
$ws.Range("H3").Value = 38.692841
$ws.Range("I3").Value = 0.1295258291743358
$ws.Range("J3").Value = 0.1295258291743358
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.307744
$ws.Range("N3").Value = 0.9232319999999999
$ws.Range("O3").Value = 0.0007278787860563825
$ws.Range("P3").Value = 0.0007278787860563825
$ws.Range("Q3").Value = 3.969163220234666
$ws.Range("R3").Value = 35.722468982112
$ws.Range("S3").Value = [double]"9.427910330236191E-05"
$ws.Range("T3").Value = [double]"9.427910330236191E-05"

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mif"
$ws.Range("C4").Value = "Cxcr4"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.89761366666667
$ws.Range("H4").Value = 38.692841
$ws.Range("I4").Value = 0.1295258291743358
$ws.Range("J4").Value = 0.1295258291743358
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 225.0171966666667
$ws.Range("N4").Value = 675.05159
$ws.Range("O4").Value = 0.5322126311204886
$ws.Range("P4").Value = 0.5322126311204886
$ws.Range("Q4").Value = 2902.184870963021
$ws.Range("R4").Value = 26119.66383866719
$ws.Range("S4").Value = 0.06893528234293619
$ws.Range("T4").Value = 0.06893528234293619

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Mif"
$ws.Range("C5").Value = "Cxcr4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.89761366666667
$ws.Range("H5").Value = 38.692841
$ws.Range("I5").Value = 0.1295258291743358
$ws.Range("J5").Value = 0.1295258291743358
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 35.56644499999999
$ws.Range("N5").Value = 106.699335
$ws.Range("O5").Value = 0.08412206512861695
$ws.Range("P5").Value = 0.08412206512861695
$ws.Range("Q5").Value = 458.7222671067482
$ws.Range("R5").Value = 4128.500403960735
$ws.Range("S5").Value = 0.01089598023764159
$ws.Range("T5").Value = 0.01089598023764159

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mif"
$ws.Range("C6").Value = "Cxcr4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.69923533333333
$ws.Range("H6").Value = 53.097706
$ws.Range("I6").Value = 0.1777466895466555
$ws.Range("J6").Value = 0.1777466895466556
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 161.9042856666667
$ws.Range("N6").Value = 485.712857
$ws.Range("O6").Value = 0.3829374249648381
$ws.Range("P6").Value = 0.3829374249648381
$ws.Range("Q6").Value = 2865.582053489561
$ws.Range("R6").Value = 25790.23848140605
$ws.Range("S6").Value = 0.06806585959102077
$ws.Range("T6").Value = 0.06806585959102078

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mif"
$ws.Range("C7").Value = "Cxcr4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.69923533333333
$ws.Range("H7").Value = 53.097706
$ws.Range("I7").Value = 0.1777466895466555
$ws.Range("J7").Value = 0.1777466895466556
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.307744
$ws.Range("N7").Value = 0.9232319999999999
$ws.Range("O7").Value = 0.0007278787860563825
$ws.Range("P7").Value = 0.0007278787860563825
$ws.Range("Q7").Value = 5.446833478421333
$ws.Range("R7").Value = 49.021501305792
$ws.Range("S7").Value = 0.0001293780446127603
$ws.Range("T7").Value = 0.0001293780446127604

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Mif"
$ws.Range("C8").Value = "Cxcr4"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 17.69923533333333
$ws.Range("H8").Value = 53.097706
$ws.Range("I8").Value = 0.1777466895466555
$ws.Range("J8").Value = 0.1777466895466556
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 225.0171966666667
$ws.Range("N8").Value = 675.05159
$ws.Range("O8").Value = 0.5322126311204886
$ws.Range("P8").Value = 0.5322126311204886
$ws.Range("Q8").Value = 3982.632317850283
$ws.Range("R8").Value = 35843.69086065255
$ws.Range("S8").Value = 0.09459903331658219
$ws.Range("T8").Value = 0.0945990333165822

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Mif"
$ws.Range("C9").Value = "Cxcr4"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 17.69923533333333
$ws.Range("H9").Value = 53.097706
$ws.Range("I9").Value = 0.1777466895466555
$ws.Range("J9").Value = 0.1777466895466556
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 35.56644499999999
$ws.Range("N9").Value = 106.699335
$ws.Range("O9").Value = 0.08412206512861695
$ws.Range("P9").Value = 0.08412206512861695
$ws.Range("Q9").Value = 629.4988800250566
$ws.Range("R9").Value = 5665.48992022551
$ws.Range("S9").Value = 0.01495241859443981
$ws.Range("T9").Value = 0.01495241859443982

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Mif"
$ws.Range("C10").Value = "Cxcr4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 30.16920766666667
$ws.Range("H10").Value = 90.507623
$ws.Range("I10").Value = 0.302977879439589
$ws.Range("J10").Value = 0.302977879439589
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 161.9042856666667
$ws.Range("N10").Value = 485.712857
$ws.Range("O10").Value = 0.3829374249648381
$ws.Range("P10").Value = 0.3829374249648381
$ws.Range("Q10").Value = 4884.52401640099
$ws.Range("R10").Value = 43960.71614760891
$ws.Range("S10").Value = 0.1160215689739034
$ws.Range("T10").Value = 0.1160215689739034

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Mif"
$ws.Range("C11").Value = "Cxcr4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 30.16920766666667
$ws.Range("H11").Value = 90.507623
$ws.Range("I11").Value = 0.302977879439589
$ws.Range("J11").Value = 0.302977879439589
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.307744
$ws.Range("N11").Value = 0.9232319999999999
$ws.Range("O11").Value = 0.0007278787860563825
$ws.Range("P11").Value = 0.0007278787860563825
$ws.Range("Q11").Value = 9.284392644170666
$ws.Range("R11").Value = 83.55953379753599
$ws.Range("S11").Value = 0.000220531171088425
$ws.Range("T11").Value = 0.000220531171088425

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Mif"
$ws.Range("C12").Value = "Cxcr4"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 30.16920766666667
$ws.Range("H12").Value = 90.507623
$ws.Range("I12").Value = 0.302977879439589
$ws.Range("J12").Value = 0.302977879439589
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 225.0171966666667
$ws.Range("N12").Value = 675.05159
$ws.Range("O12").Value = 0.5322126311204886
$ws.Range("P12").Value = 0.5322126311204886
$ws.Range("Q12").Value = 6788.590534807841
$ws.Range("R12").Value = 61097.31481327057
$ws.Range("S12").Value = 0.1612486543878498
$ws.Range("T12").Value = 0.1612486543878498

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Mif"
$ws.Range("C13").Value = "Cxcr4"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 30.16920766666667
$ws.Range("H13").Value = 90.507623
$ws.Range("I13").Value = 0.302977879439589
$ws.Range("J13").Value = 0.302977879439589
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 35.56644499999999
$ws.Range("N13").Value = 106.699335
$ws.Range("O13").Value = 0.08412206512861695
$ws.Range("P13").Value = 0.08412206512861695
$ws.Range("Q13").Value = 1073.011465170078
$ws.Range("R13").Value = 9657.103186530703
$ws.Range("S13").Value = 0.02548712490674736
$ws.Range("T13").Value = 0.02548712490674736

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Mif"
$ws.Range("C14").Value = "Cxcr4"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 38.809555
$ws.Range("H14").Value = 116.428665
$ws.Range("I14").Value = 0.3897496018394196
$ws.Range("J14").Value = 0.3897496018394196
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 161.9042856666667
$ws.Range("N14").Value = 485.712857
$ws.Range("O14").Value = 0.3829374249648381
$ws.Range("P14").Value = 0.3829374249648381
$ws.Range("Q14").Value = 6283.433279316213
$ws.Range("R14").Value = 56550.89951384591
$ws.Range("S14").Value = 0.1492497089094582
$ws.Range("T14").Value = 0.1492497089094583

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Mif"
$ws.Range("C15").Value = "Cxcr4"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 38.809555
$ws.Range("H15").Value = 116.428665
$ws.Range("I15").Value = 0.3897496018394196
$ws.Range("J15").Value = 0.3897496018394196
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.307744
$ws.Range("N15").Value = 0.9232319999999999
$ws.Range("O15").Value = 0.0007278787860563825
$ws.Range("P15").Value = 0.0007278787860563825
$ws.Range("Q15").Value = 11.94340769392
$ws.Range("R15").Value = 107.49066924528
$ws.Range("S15").Value = 0.0002836904670528352
$ws.Range("T15").Value = 0.0002836904670528352

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Mif"
$ws.Range("C16").Value = "Cxcr4"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 38.809555
$ws.Range("H16").Value = 116.428665
$ws.Range("I16").Value = 0.3897496018394196
$ws.Range("J16").Value = 0.3897496018394196
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 225.0171966666667
$ws.Range("N16").Value = 675.05159
$ws.Range("O16").Value = 0.5322126311204886
$ws.Range("P16").Value = 0.5322126311204886
$ws.Range("Q16").Value = 8732.817269980818
$ws.Range("R16").Value = 78595.35542982737
$ws.Range("S16").Value = 0.2074296610731203
$ws.Range("T16").Value = 0.2074296610731203

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Mif"
$ws.Range("C17").Value = "Cxcr4"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 38.809555
$ws.Range("H17").Value = 116.428665
$ws.Range("I17").Value = 0.3897496018394196
$ws.Range("J17").Value = 0.3897496018394196
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 35.56644499999999
$ws.Range("N17").Value = 106.699335
$ws.Range("O17").Value = 0.08412206512861695
$ws.Range("P17").Value = 0.08412206512861695
$ws.Range("Q17").Value = 1380.317903381975
$ws.Range("R17").Value = 12422.86113043778
$ws.Range("S17").Value = 0.03278654138978818
$ws.Range("T17").Value = 0.03278654138978818
